$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 148.9466983333333
$ws.Range("H2").Value = 446.840095
$ws.Range("I2").Value = 0.5255511750002251
$ws.Range("J2").Value = 0.5255511750002251
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 169.629438
$ws.Range("N2").Value = 508.888314
$ws.Range("O2").Value = 0.7428377317484701
$ws.Range("P2").Value = 0.7428377317484702
$ws.Range("Q2").Value = 25265.74473023887
$ws.Range("R2").Value = 227391.7025721498
$ws.Range("S2").Value = 0.3903992427549104
$ws.Range("T2").Value = 0.3903992427549106

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 148.9466983333333
$ws.Range("H3").Value = 446.840095
$ws.Range("I3").Value = 0.5255511750002251
$ws.Range("J3").Value = 0.5255511750002251
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.9848756666666668
$ws.Range("N3").Value = 2.954627
$ws.Range("O3").Value = 0.004312947180081616
$ws.Range("P3").Value = 0.004312947180081616
$ws.Range("Q3").Value = 146.6939788188406
$ws.Range("R3").Value = 1320.245809369565
$ws.Range("S3").Value = 0.0022666744582058
$ws.Range("T3").Value = 0.0022666744582058

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 148.9466983333333
$ws.Range("H4").Value = 446.840095
$ws.Range("I4").Value = 0.5255511750002251
$ws.Range("J4").Value = 0.5255511750002251
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 54.620752
$ws.Range("N4").Value = 163.862256
$ws.Range("O4").Value = 0.2391940691454494
$ws.Range("P4").Value = 0.2391940691454494
$ws.Range("Q4").Value = 8135.580670883814
$ws.Range("R4").Value = 73220.22603795433
$ws.Range("S4").Value = 0.125708724092476
$ws.Range("T4").Value = 0.125708724092476

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 148.9466983333333
$ws.Range("H5").Value = 446.840095
$ws.Range("I5").Value = 0.5255511750002251
$ws.Range("J5").Value = 0.5255511750002251
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.118221666666667
$ws.Range("N5").Value = 9.354665000000001
$ws.Range("O5").Value = 0.01365525192599884
$ws.Range("P5").Value = 0.01365525192599884
$ws.Range("Q5").Value = 464.4488219214639
$ws.Range("R5").Value = 4180.039397293175
$ws.Range("S5").Value = 0.007176533694632778
$ws.Range("T5").Value = 0.007176533694632778

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 103.6567713333333
$ws.Range("H6").Value = 310.970314
$ws.Range("I6").Value = 0.3657478720948015
$ws.Range("J6").Value = 0.3657478720948015
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 169.629438
$ws.Range("N6").Value = 508.888314
$ws.Range("O6").Value = 0.7428377317484701
$ws.Range("P6").Value = 0.7428377317484702
$ws.Range("Q6").Value = 17583.23986616785
$ws.Range("R6").Value = 158249.1587955106
$ws.Range("S6").Value = 0.2716913196987319
$ws.Range("T6").Value = 0.271691319698732

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 103.6567713333333
$ws.Range("H7").Value = 310.970314
$ws.Range("I7").Value = 0.3657478720948015
$ws.Range("J7").Value = 0.3657478720948015
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.9848756666666668
$ws.Range("N7").Value = 2.954627
$ws.Range("O7").Value = 0.004312947180081616
$ws.Range("P7").Value = 0.004312947180081616
$ws.Range("Q7").Value = 102.0890317714309
$ws.Range("R7").Value = 918.8012859428782
$ws.Range("S7").Value = 0.001577451253572126
$ws.Range("T7").Value = 0.001577451253572126

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 103.6567713333333
$ws.Range("H8").Value = 310.970314
$ws.Range("I8").Value = 0.3657478720948015
$ws.Range("J8").Value = 0.3657478720948015
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 54.620752
$ws.Range("N8").Value = 163.862256
$ws.Range("O8").Value = 0.2391940691454494
$ws.Range("P8").Value = 0.2391940691454494
$ws.Range("Q8").Value = 5661.81080011871
$ws.Range("R8").Value = 50956.29720106839
$ws.Range("S8").Value = 0.08748472180764494
$ws.Range("T8").Value = 0.08748472180764494

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 103.6567713333333
$ws.Range("H9").Value = 310.970314
$ws.Range("I9").Value = 0.3657478720948015
$ws.Range("J9").Value = 0.3657478720948015
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.118221666666667
$ws.Range("N9").Value = 9.354665000000001
$ws.Range("O9").Value = 0.01365525192599884
$ws.Range("P9").Value = 0.01365525192599884
$ws.Range("Q9").Value = 323.2247902683123
$ws.Range("R9").Value = 2909.02311241481
$ws.Range("S9").Value = 0.004994379334852517
$ws.Range("T9").Value = 0.004994379334852517

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.087957333333333
$ws.Range("H10").Value = 3.263872
$ws.Range("I10").Value = 0.003838804493697762
$ws.Range("J10").Value = 0.003838804493697762
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 169.629438
$ws.Range("N10").Value = 508.888314
$ws.Range("O10").Value = 0.7428377317484701
$ws.Range("P10").Value = 0.7428377317484702
$ws.Range("Q10").Value = 184.549591021312
$ws.Range("R10").Value = 1660.946319191808
$ws.Range("S10").Value = 0.00285160882272428
$ws.Range("T10").Value = 0.00285160882272428

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.087957333333333
$ws.Range("H11").Value = 3.263872
$ws.Range("I11").Value = 0.003838804493697762
$ws.Range("J11").Value = 0.003838804493697762
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.9848756666666668
$ws.Range("N11").Value = 2.954627
$ws.Range("O11").Value = 0.004312947180081616
$ws.Range("P11").Value = 0.004312947180081616
$ws.Range("Q11").Value = 1.071502703971556
$ws.Range("R11").Value = 9.643524335744001
$ws.Range("S11").Value = 0.0000165565610159784
$ws.Range("T11").Value = 0.0000165565610159784

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.087957333333333
$ws.Range("H12").Value = 3.263872
$ws.Range("I12").Value = 0.003838804493697762
$ws.Range("J12").Value = 0.003838804493697762
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 54.620752
$ws.Range("N12").Value = 163.862256
$ws.Range("O12").Value = 0.2391940691454494
$ws.Range("P12").Value = 0.2391940691454494
$ws.Range("Q12").Value = 59.42504769058134
$ws.Range("R12").Value = 534.825429215232
$ws.Range("S12").Value = 0.0009182192675014043
$ws.Range("T12").Value = 0.0009182192675014042

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.087957333333333
$ws.Range("H13").Value = 3.263872
$ws.Range("I13").Value = 0.003838804493697762
$ws.Range("J13").Value = 0.003838804493697762
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.118221666666667
$ws.Range("N13").Value = 9.354665000000001
$ws.Range("O13").Value = 0.01365525192599884
$ws.Range("P13").Value = 0.01365525192599884
$ws.Range("Q13").Value = 3.39249212920889
$ws.Range("R13").Value = 30.53242916288
$ws.Range("S13").Value = 0.00005241984245609938
$ws.Range("T13").Value = 0.00005241984245609938

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 29.71902933333333
$ws.Range("H14").Value = 89.15708799999999
$ws.Range("I14").Value = 0.1048621484112755
$ws.Range("J14").Value = 0.1048621484112755
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 169.629438
$ws.Range("N14").Value = 508.888314
$ws.Range("O14").Value = 0.7428377317484701
$ws.Range("P14").Value = 0.7428377317484702
$ws.Range("Q14").Value = 5041.222243718847
$ws.Range("R14").Value = 45371.00019346963
$ws.Range("S14").Value = 0.07789556047210336
$ws.Range("T14").Value = 0.07789556047210337

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 29.71902933333333
$ws.Range("H15").Value = 89.15708799999999
$ws.Range("I15").Value = 0.1048621484112755
$ws.Range("J15").Value = 0.1048621484112755
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.9848756666666668
$ws.Range("N15").Value = 2.954627
$ws.Range("O15").Value = 0.004312947180081616
$ws.Range("P15").Value = 0.004312947180081616
$ws.Range("Q15").Value = 29.26954882735289
$ws.Range("R15").Value = 263.425939446176
$ws.Range("S15").Value = 0.0004522649072877107
$ws.Range("T15").Value = 0.0004522649072877107

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 29.71902933333333
$ws.Range("H16").Value = 89.15708799999999
$ws.Range("I16").Value = 0.1048621484112755
$ws.Range("J16").Value = 0.1048621484112755
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 54.620752
$ws.Range("N16").Value = 163.862256
$ws.Range("O16").Value = 0.2391940691454494
$ws.Range("P16").Value = 0.2391940691454494
$ws.Range("Q16").Value = 1623.275730896725
$ws.Range("R16").Value = 14609.48157807053
$ws.Range("S16").Value = 0.02508240397782702
$ws.Range("T16").Value = 0.02508240397782702

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 29.71902933333333
$ws.Range("H17").Value = 89.15708799999999
$ws.Range("I17").Value = 0.1048621484112755
$ws.Range("J17").Value = 0.1048621484112755
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 3.118221666666667
$ws.Range("N17").Value = 9.354665000000001
$ws.Range("O17").Value = 0.01365525192599884
$ws.Range("P17").Value = 0.01365525192599884
$ws.Range("Q17").Value = 92.67052117950222
$ws.Range("R17").Value = 834.0346906155199
$ws.Range("S17").Value = 0.001431919054057447
$ws.Range("T17").Value = 0.001431919054057447
